$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Talent_Acquisition")
$ws3 = $wb.Worksheets.Item("Employee_Details")

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS","AT","AU","AV","AW","AX","AY","AZ","BA","BB")
$headerVals = @("scenario","addressLine1","addressLine2","actualAmount","assignmentCategory","businessUnit","city","companyEntity","costCenter","country","countryCode","county","dateOfBirth","department","destinationLegalEmployer","elementName","elementNameADP","email","employmentAction","employmentActionReason","effectiveAsOfDate","erProfitCenter","firstName","gender","globalMobilityIndicator","grade","hireAction","hireDate","hireReason","hourlyPaidOrSalaried","job","lastName","legalEmployer","location","maritalStatus","nationalID","nationalIDType","newlySeparatedVeteranDischargeDate","number","password","pearsonEmailRequired","payroll","payCheckDate","personNumber","salaryAmount","salaryBasis","state","type1","type2","userName","veteranSelfIdentificationStatus","workerType","workingAsManager","zipCode")
$row2Vals = @("UPDATE_SALARY_FOR_EXISTING_EMP","527 Johns Avenue","Suite 436","'1,000.00","Full-time regular","Finance","Irving","'510","'61465","United States","United States 1","Dallas","12-Apr-1988","Finance Internal Audit - EMEA","Pearson Inc","Bilingual Indicator","ADP Auto & Home","uatnewhire67@abc.com","Global Temporary Assignment","GM Temporary Relocation","13-Feb-2019","'31",$null,"Male","International Assignment","E","Hire","28-Jan-19","Additional Hire","Salaried","Director Audit & Compliance",$null,"NCS Pearson, Inc","AR-Buenos Aires-Humboldt 1509/13","Single",$null,"Social Security Number","12-Apr-2017","2717522","Welcome123","Yes","Bi-Weekly Exempt","13-Feb-2019","3001712","80,000.00","US Annual Salary","TX","Work Phone","Home E-Mail","3259228","Not a Protected Veteran","Employee","No","75038")
$row3Vals = @("CHANGE_SALARY_BASIS_FOR_EXISTING_EMP","527 Johns Avenue","Suite 436","'1,000.00","Full-time regular","Finance","Irving","'510","'61465","United States","United States 1","Dallas","12-Apr-1988","Finance Internal Audit - EMEA","Pearson Inc","Bilingual Indicator","ADP Auto & Home","uatnewhire67@abc.com","Global Temporary Assignment","GM Temporary Relocation","13-Feb-2019","'31",$null,"Male","International Assignment","E","Hire","28-Jan-19","Additional Hire","Salaried","Director Audit & Compliance",$null,"NCS Pearson, Inc","AR-Buenos Aires-Humboldt 1509/13","Single",$null,"Social Security Number","12-Apr-2017","2717522","Welcome123","Yes","Bi-Weekly Exempt","13-Feb-2019","3259949","80,000.00","US Annual Salary","TX","Work Phone","Home E-Mail","3040301","Not a Protected Veteran","Employee","No","75038")
$row4Vals = @("EDIT_SALARY_PROPOSAL_REASON","527 Johns Avenue","Suite 436","'1,000.00","Full-time regular","Finance","Irving","'510","'61465","United States","United States 1","Dallas","12-Apr-1988","Finance Internal Audit - EMEA","Pearson Inc","Bilingual Indicator","ADP Auto & Home","uatnewhire67@abc.com","Global Temporary Assignment","GM Temporary Relocation","13-Feb-2019","'31",$null,"Male","International Assignment","E","Hire","28-Jan-19","Additional Hire","Salaried","Director Audit & Compliance",$null,"NCS Pearson, Inc","AR-Buenos Aires-Humboldt 1509/13","Single",$null,"Social Security Number","12-Apr-2017","2717522","Welcome123","Yes","Bi-Weekly Exempt","13-Feb-2019","3259949","15.00","US Annual Salary","TX","Work Phone","Home E-Mail","3259228","Not a Protected Veteran","Employee","No","75038")
$dataStyles = @("0","1","1","2","1","1","1","2","2","1","1","1","1","1","1","1","1","1","1","1","1","2","1","1","1","1","1","1","1","1","1","1","1","1","1","1","1","1","1","1","1","1","1","1","1","1","1","1","1","1","1","1","1","1")

# --- Header row (row 1): plain values, default (unstyled) cells ---
for ($i = 0; $i -lt $cols.Length; $i++) {
    $c = $cols[$i]
    $ws3.Range("$c`1").Value = $headerVals[$i]
}

# --- Data rows 2-4 (copied from Talent_Acquisition rows 7,8,9) ---
$dataRows = @(2,3,4)
$dataVals = @($row2Vals, $row3Vals, $row4Vals)

for ($r = 0; $r -lt $dataRows.Length; $r++) {
    $rowNum = $dataRows[$r]
    $vals = $dataVals[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $c = $cols[$i]
        $style = $dataStyles[$i]
        $cell = $ws3.Range("$c$rowNum")
        if ($style -eq "1" -or $style -eq "2") {
            $cell.NumberFormat = "@"
        }
        if ($vals[$i] -ne $null) {
            $cell.Value = $vals[$i]
        }
    }
}

# --- Comment on AR2 (mirrors the threaded comment on Talent_Acquisition!AR7) ---
$commentText = "[Threaded comment]" + [char]10 + "Your version of Excel allows you to read this threaded comment; however, any edits to it will get removed if the file is opened in a newer version of Excel. Learn more: https://go.microsoft.com/fwlink/?linkid=870924" + [char]10 + "Comment:" + [char]10 + "    3001660"
$ws3.Range("AR2").AddComment($commentText)

# --- View / selection state ---
$ws1.Activate()
$ws1.Range("A1:XFD1048576").Select()

$ws3.Activate()
$ws3.Range("A8").Select()
